# Auto-generated: re-apply the updated crypto price/volume figures
# captured by the GitHub Actions scraper commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.282.66"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "1.865.70"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.66%  "
$ws.Range("D5").Value = "'239.88"
$ws.Range("E5").Value = "  +3.53%  "
$ws.Range("D6").Value = "'0.625"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").Value = "'42.62"
$ws.Range("E8").Value = "  +6.63%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "2.132.77"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "'11.50"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "1.864.12"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "'4.74"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "35.253.59"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "'70.06"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "0.0₃0798"
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").Value = "'241.61"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "'12.27"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").Value = "'4.76"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "'169.56"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").Value = "'1.91"
$ws.Range("E26").Value = "  +25.08%  "
$ws.Range("E27").Value = "  +4.22%  "
$ws.Range("D28").Value = "'17.76"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "'0.0564"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").Value = "'1.83"
$ws.Range("E33").Value = "  +27.90%  "
$ws.Range("D34").Value = "'4.05"
$ws.Range("E34").Value = "  +2.21%  "
$ws.Range("E35").Value = "  +8.20%  "
$ws.Range("D36").Value = "'0.823"
$ws.Range("E36").Value = "  +18.09%  "
$ws.Range("D37").Value = "'1.31"
$ws.Range("E37").Value = "  +6.75%  "
$ws.Range("E38").Value = "  +3.53%  "
$ws.Range("E39").Value = "  +4.58%  "
$ws.Range("D40").Value = "'90.55"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").Value = "1.347.16"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("D42").Value = "'15.29"
$ws.Range("D43").Value = "'0.0603"
$ws.Range("E43").Value = "  +15.49%  "
$ws.Range("D44").Value = "'2.35"
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "'12.42"
$ws.Range("E46").Value = "  +44.19%  "
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").Value = "'6.61"
$ws.Range("E48").Value = "  +5.00%  "
$ws.Range("D49").Value = "2.049.18"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").Value = "'0.0688"
$ws.Range("E50").Value = "  +3.31%  "
$ws.Range("D51").Value = "'3.43"
$ws.Range("E51").Value = "  +1.52%  "
